# cambios de agosto, puntos fe de ratas e historico
#
# Moves the reported period from Q1 2022 (Jan-Mar) to Q2 2022 (Apr-Jun),
# updates the validation/update dates, refreshes a couple of normativity
# hyperlink URLs (http -> https, and the Hidalgo official-gazette link to
# the newer "alcance-9" edition), drops the now-stale hyperlink objects,
# shrinks the Hidden_13 list-validation range on column D, and moves the
# active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

# --- Row 8: Decreto de creación ---------------------------------------
$ws.Range("B8").Value = 44652
$ws.Range("C8").Value = 44742
$ws.Range("H8").Value = "https://www.upp.edu.mx/normatividad/files/interna/decretos/decreto-de-creacion-vigente-04_02_2008.pdf"
$ws.Range("J8").Value = 44753
$ws.Range("K8").Value = 44753

# --- Row 9: Estatuto Orgánico -------------------------------------------
$ws.Range("B9").Value = 44652
$ws.Range("C9").Value = 44742
$ws.Range("H9").Value = "https://www.upp.edu.mx/normatividad/files/interna/estatutos/estatuto-organico-de-la-universidad-politecnica-de-pachuca-25_08_2008.pdf"
$ws.Range("J9").Value = 44753
$ws.Range("K9").Value = 44753

# --- Row 10: Cuotas y Tarifas --------------------------------------------
$ws.Range("B10").Value = 44652
$ws.Range("C10").Value = 44742
$ws.Range("F10").Value = 44539
$ws.Range("G10").Value = 44561
$ws.Range("H10").Value = "https://periodico.hidalgo.gob.mx/?tribe_events=periodico-oficial-alcance-9-del-31-de-diciembre-de-2021"
$ws.Range("J10").Value = 44753
$ws.Range("K10").Value = 44753

# The source documents no longer carry live hyperlink objects.
$ws.Range("H8").Hyperlinks.Delete()
$ws.Range("H9").Hyperlinks.Delete()
$ws.Range("H10").Hyperlinks.Delete()

# Shrink the Hidden_13 list-validation range applied to column D.
$ws.Range("D11:D200").Validation.Delete()
$v = $ws.Range("D11:D107").Validation
$v.Add(3, 1, 1, "=Hidden_13")
$v.IgnoreBlank = $true
$v.InCellDropdown = $true
$v.ShowInput = $false
$v.ShowError = $true

# Move the active selection/scroll position.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
[void]$ws.Range("C15").Select()
